# StoryBoardSOS.pptx edit script
#
# 1. Update the cached "datetimeFigureOut" footer date from 3/21/2016 to
#    3/28/2016 on the slide master and all 11 slide layouts.
# 2. Slide 2 ("Sent" state storyboard), inside the "Group 34" shape:
#      - merge the two separate text runs "SOS" + "! I need help!" in the
#        "Rectangle 21" shape into a single run "SOS! I need help!".
#      - recolor the "Oval 4" shape's fill and outline from the hard-coded
#        red (FF0000) to the theme color accent2.

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# 1. Date placeholder text (slide master + every custom layout)
# ---------------------------------------------------------------------
function Update-DatePlaceholder($shapes, $newText) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shp = $shapes.Item($i)
        if ($shp.Name -like "Date Placeholder*") {
            $shp.TextFrame.TextRange.Text = $newText
        }
    }
}

$master = $p.SlideMaster
Update-DatePlaceholder $master.Shapes "3/28/2016"

for ($L = 1; $L -le $master.CustomLayouts.Count; $L++) {
    $layout = $master.CustomLayouts.Item($L)
    Update-DatePlaceholder $layout.Shapes "3/28/2016"
}

# ---------------------------------------------------------------------
# 2. Slide 2 edits
# ---------------------------------------------------------------------
$slide2 = $p.Slides.Item(2)
$group34 = $slide2.Shapes.Item(1)

# -- Oval 4: red -> accent2 theme color (fill + line) -----------------
$oval = $group34.GroupItems.Item(2)
if ($oval.Name -eq "Oval 4") {
    $oval.Fill.ForeColor.ObjectThemeColor = 6   # msoThemeColorAccent2
    $oval.Line.ForeColor.ObjectThemeColor = 6   # msoThemeColorAccent2
}

# -- Rectangle 21: merge "SOS" + "! I need help!" into one run --------
$rect = $group34.GroupItems.Item(3)
if ($rect.Name -eq "Rectangle 21") {
    $para1 = $rect.TextFrame.TextRange.Paragraphs(1)
    # Force a real rebuild of the paragraph's runs: re-assigning the
    # exact same text is treated as a no-op by the engine, so bounce
    # through a temporary value first.
    $para1.Text = "SOS! I need help!__tmp__"
    $para1.Text = "SOS! I need help!"
}
